$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update handoff/handback datetimes for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 14:46:07"
$wsZhCn.Range("H2").Value = "2016-03-11 14:46:37"

# "de-de" sheet: update handoff/handback datetimes for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 14:46:11"
$wsDeDe.Range("H2").Value = "2016-03-11 14:46:43"
